# Addiction DictionaryBase class and Spanish Terms test
#
# 1. Select A1:B2 on the "DictionaryPopup" sheet (2nd sheet) *before* adding
#    the new sheet, so that the newly added sheet ends up being the final
#    active / tab-selected sheet (matches workbookView activeTab="2").
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A1:B2").Select() | Out-Null

# 2. Add the new "TermsSpanish" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "TermsSpanish"

# 3. Populate the new sheet's values -- the column-B (ContentType) values for
#    the new rows are entered before the column-A (Path) values so the
#    shared-string table ends up in the same append order as the target
#    workbook.
$newSheet.Range("A1").Value = "Path"
$newSheet.Range("B1").Value = "ContentType"

$newSheet.Range("A2").Value = "/espanol/publicaciones/diccionario"
$newSheet.Range("B2").Value = "Dictionary Page"

$newSheet.Range("B3").Value = "Dictionary Search Page"
$newSheet.Range("B4").Value = "Dictionary Page Expand"

$newSheet.Range("A3").Value = "/espanol/publicaciones/diccionario/buscar"
$newSheet.Range("A4").Value = "/espanol/publicaciones/diccionario?expand=D"

# 4. Copy the header formatting (bold font + gray fill) from the existing
#    "DictionaryPage" header row so the new header reuses the same style.
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Match the column widths used by the other dictionary sheets as closely
#    as possible.
$newSheet.Columns.Item(1).ColumnWidth = 42.592447916666664
$newSheet.Columns.Item(2).ColumnWidth = 21.166666666666668

# 6. Leave the selection on the new sheet parked one row below the data,
#    mirroring the other sheets' "next empty row" selection convention.
$newSheet.Range("A5").Select() | Out-Null
